$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill the "duplicate_image_filename" column (E) with "NA" for every
# data row in the practice (rows 2-5) and trial (rows 6-21) blocks.
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
